$d = $word.ActiveDocument

# 1) Merge the two "C = 1," runs (" = " + "0.8639") into a single run.
#    Word normalizes run boundaries automatically when Find/Replace
#    rewrites text that spans multiple runs with identical formatting.
$d.Content.Find.Execute(" = 0.8639", $false, $false, $false, $false, $false, `
                         $true, 1, $false, " = 0.8639", 2)

# 2) The "C = 0.1" paragraph has the exact same " = " + "0.8639" run split;
#    running the same replace again (Find continues over the whole story)
#    merges that occurrence too.
$d.Content.Find.Execute(" = 0.8639", $false, $false, $false, $false, $false, `
                         $true, 1, $false, " = 0.8639", 2)

# 3) Merge the many small runs making up the classifier equation into one
#    run (they all share the same rPr/style, so a full-text replace collapses
#    them into a single run).
$eq = "-0.08158492 -0.0010065348*X1 -0.0011729048*X2 -0.0016261967*X3 + 0.0030064203*X8+ 1.0049405641*X9 -0.0028259432*X10 + 0.0002600295*X11 -0.0005349551*X12 -0.0012283758*X14 + 0.106363399*X15 = 0"
$d.Content.Find.Execute($eq, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $eq, 2)

# 4) Replace "How to remove itself in KNN. " with the new Q3 cross-validation
#    note, and push the trailing _GoBack bookmark into its own paragraph by
#    inserting a paragraph break (^p) right after the new text.
$d.Content.Find.Execute("How to remove itself in KNN. ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "K = 10, accuracy = 0.924^p", 2)

# Locate the now-isolated bookmark paragraph (the last paragraph in the
# document) and append the new "Question 3." and blank-space paragraphs
# after it, so the bookmark keeps living alone in its own paragraph.
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkPara.Range.InsertParagraphAfter()

$q3Para = $d.Paragraphs.Item($d.Paragraphs.Count)
$q3Para.Range.Text = "Question 3."

$q3Para = $d.Paragraphs.Item($d.Paragraphs.Count)
$q3Para.Range.InsertParagraphAfter()
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Text = " "
